$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the CreatedBy name (C2): Peter Styliadis -> Su Chee Tay
$ws.Range("C2").Value = "Su Chee Tay"

# Update the Email (F2): Peter@myemail.com -> SuChee@myemail.com, and turn it into
# a mailto: hyperlink (this also applies the built-in "Hyperlink" style/font).
$ws.Range("F2").Value = "SuChee@myemail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "mailto:SuChee@myemail.com")

# The welcome/description row no longer needs as much vertical space.
$ws.Rows.Item(2).RowHeight = 75

# Move the active selection/cursor to A4.
$ws.Range("A4").Select()
